# fix list student in grade default file
# Adds a missing student row (StudentId + Full name) to "sheet 1".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New student record appended as row 5, matching the formatting of the
# existing data rows (A3:B4): numeric id left-aligned, name as plain text.
$ws.Range("A5").Value = 20120600
$ws.Range("B5").Value = "Nguyển Văn A"
$ws.Range("A5").HorizontalAlignment = -4131

# Leave the selection where the user ended up after entering the new row.
$ws.Range("D3").Select()
